$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.216.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "'1.600.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'52.05"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.39%  "
$ws.Range("D9").Value = "'0.3638"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").Value = "'1.270"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "'0.08124"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").Value = "'22.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").Value = "'6.572"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "'7.408"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "'0.00001246"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").Value = "'1.602.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "'94.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.45%  "
$ws.Range("D19").Value = "'0.06925"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").Value = "'18.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "'6.526"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'12.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").Value = "'23.231.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.448"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.83%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'3.043"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.40%  "
$ws.Range("D27").Value = "'21.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("D28").Value = "'149.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").Value = "'134.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").Value = "'2.394"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.88%  "
$ws.Range("D32").Value = "'6.692"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("D33").Value = "'1.781.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D34").Value = "'0.9595"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").Value = "'0.07477"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("D36").Value = "'10.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("D37").Value = "'0.02744"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("D38").Value = "'0.2533"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").Value = "'0.08801"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").Value = "'6.074"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").Value = "'1.386"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("D42").Value = "'0.7084"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").Value = "'12.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "'15.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("D45").Value = "'0.6528"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'2.313"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.21%  "
$ws.Range("D48").Value = "'4.011"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("D49").Value = "'132.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("D50").Value = "'0.07925"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "'1.201"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.23%  "
